$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Price column (D) updates - forced to text to preserve exact formatting
Set-TextValue "D2" "69.266.85"
Set-TextValue "D3" "3.739.91"
Set-TextValue "D5" "614.18"
Set-TextValue "D6" "187.80"
Set-TextValue "D7" "3.734.83"
Set-TextValue "D8" "0.642"
Set-TextValue "D10" "0.723"
Set-TextValue "D11" "0.163"
Set-TextValue "D12" "57.43"
Set-TextValue "D13" "0.0000294"
Set-TextValue "D14" "10.71"
Set-TextValue "D15" "4.329.27"
Set-TextValue "D16" "3.736.57"
Set-TextValue "D17" "19.44"
Set-TextValue "D18" "13.12"
Set-TextValue "D19" "1.14"
Set-TextValue "D20" "0.127"
Set-TextValue "D21" "69.044.76"
Set-TextValue "D22" "414.55"
Set-TextValue "D23" "4.65"
Set-TextValue "D24" "89.69"
Set-TextValue "D25" "3.07"
Set-TextValue "D26" "13.00"
Set-TextValue "D27" "11.12"
Set-TextValue "D29" "3.81"
Set-TextValue "D30" "9.73"
Set-TextValue "D31" "33.45"
Set-TextValue "D32" "7.41"
Set-TextValue "D33" "12.82"
Set-TextValue "D35" "628.13"
Set-TextValue "D36" "45.11"
Set-TextValue "D37" "66.32"
Set-TextValue "D38" "0.0₃0845"
Set-TextValue "D39" "0.422"
Set-TextValue "D41" "0.999"
Set-TextValue "D42" "0.141"
Set-TextValue "D43" "3.08"
Set-TextValue "D44" "0.0448"
Set-TextValue "D45" "2.66"
Set-TextValue "D47" "2.856.85"
Set-TextValue "D48" "2.78"
Set-TextValue "D49" "9.21"
Set-TextValue "D50" "2.74"
Set-TextValue "D51" "3.14"

# Coin name / Link / Volume column updates
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  +5.19%  "
$ws.Range("E6").Value = "  +6.64%  "
$ws.Range("E7").Value = "  +0.92%  "
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("E12").Value = "  +8.24%  "
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("E15").Value = "  +0.54%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("B19").Value = "Polygon"
$ws.Range("C19").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  +1.03%  "
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("E30").Value = "  +1.31%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  -8.51%  "
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("E35").Value = "  +4.81%  "
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -9.77%  "
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("E46").Value = "  +4.04%  "
$ws.Range("E47").Value = "  +3.76%  "
$ws.Range("E48").Value = "  +5.21%  "
$ws.Range("E49").Value = "  -3.94%  "
$ws.Range("E50").Value = "  -19.26%  "
$ws.Range("E51").Value = "  +0.23%  "
